$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns (D and E) before the existing "TAILLE" column,
# which shifts it from D to F.
$ws.Columns.Item(4).Insert()
$ws.Columns.Item(4).Insert()

# Fill the new columns with their values, row by row.
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 4).Value2 = 10
    $ws.Cells.Item($r, 5).Value2 = 11
}

# Match the width of the neighbouring data columns.
$ws.Columns.Item(4).ColumnWidth = $ws.Columns.Item(2).ColumnWidth
$ws.Columns.Item(5).ColumnWidth = $ws.Columns.Item(2).ColumnWidth

# Adjust row heights.
$ws.Rows.Item(1).RowHeight = 46
for ($r = 2; $r -le 25; $r++) {
    $ws.Rows.Item($r).RowHeight = 23
}

$ws.Range("D2:E25").Select()
